# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for rows 2-19, replacing the previous Strike# values.
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
